$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows.Item(442).Insert()

$ws.Range("A442").Value = 3
$ws.Range("B442").Value = "Femacal de La Calera"
$ws.Range("C442").Value = "Coquimbo"
$ws.Range("D442").Value = 45142
$ws.Range("E442").Value = 5
$ws.Range("F442").Value = 100112001
$ws.Range("G442").Value = "Berenjena"
$ws.Range("H442").Value = "Sin especificar"
$ws.Range("I442").Value = "Primera"
$ws.Range("J442").Value = 55
$ws.Range("K442").Value = 7500
$ws.Range("L442").Value = 7500
$ws.Range("M442").Value = 7500
$ws.Range("N442").Value = '$/caja 60 unidades'
$ws.Range("O442").Value = "Región de Arica y Parinacota"
$ws.Range("P442").Value = 125
$ws.Range("Q442").Value = 60
$ws.Range("R442").Value = "Hortaliza"
